$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 - reorder "Recorded By" list
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 3 - reorder "Recorded By" list
$ws.Range("G3").Value = "System, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# Row 4 - add recorders and update attendance count
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H4").Value = "48/251"

# Row 9 - reorder "Recorded By" list
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Updated Average Attendance % statistics (force text so Excel doesn't
# reinterpret the percent string as a numeric percentage value)
$ws.Range("L10").Value = "'24.3%"
$ws.Range("S15").Value = "'24.3%"

# Row 28 - reorder "Recorded By" list
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
